$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value2 = 4293.0835
$ws.Range("I40").Value2 = 3247.7273
$ws.Range("K40").Value2 = 3247.7273
$ws.Range("M40").Value2 = -3072.7273
$ws.Range("H51").Value2 = 5000
$ws.Range("J51").Value2 = 5000
$ws.Range("L51").Value2 = 5000
$ws.Range("N51").Value2 = -5968
$ws.Range("H69").Value2 = 33414224
$ws.Range("I69").Value2 = 19999
$ws.Range("J69").Value2 = 37124692
$ws.Range("K69").Value2 = 59997
$ws.Range("L69").Value2 = 111374076
$ws.Range("N69").Value2 = -111375824
$ws.Range("M69").Value2 = -59123
$ws.Range("H70").Value2 = 4845.769
$ws.Range("I70").Value2 = 3000
$ws.Range("J70").Value2 = 4999.5835
$ws.Range("K70").Value2 = 9000
$ws.Range("L70").Value2 = 14998.7505
$ws.Range("N70").Value2 = -15538.7505
$ws.Range("M70").Value2 = -8730
$ws.Range("H72").Value2 = 33414224
$ws.Range("I72").Value2 = 19999
$ws.Range("J72").Value2 = 37124692
$ws.Range("K72").Value2 = 179991
$ws.Range("L72").Value2 = 334122228
$ws.Range("N72").Value2 = -334130964
$ws.Range("M72").Value2 = -175623
$ws.Range("H73").Value2 = 4845.769
$ws.Range("I73").Value2 = 3000
$ws.Range("J73").Value2 = 4999.5835
$ws.Range("K73").Value2 = 9000
$ws.Range("L73").Value2 = 14998.7505
$ws.Range("N73").Value2 = -16870.7505
$ws.Range("M73").Value2 = -8064
$ws.Range("H76").Value2 = 3224.25
$ws.Range("I76").Value2 = 999
$ws.Range("J76").Value2 = 9900
$ws.Range("K76").Value2 = 999
$ws.Range("L76").Value2 = 9900
$ws.Range("M76").Value2 = -684
$ws.Range("N76").Value2 = -10530
$ws.Range("H79").Value2 = 3224.25
$ws.Range("I79").Value2 = 999
$ws.Range("J79").Value2 = 9900
$ws.Range("K79").Value2 = 999
$ws.Range("L79").Value2 = 9900
$ws.Range("M79").Value2 = 93
$ws.Range("N79").Value2 = -12084
$ws.Range("H94").Value2 = 1168.9333
$ws.Range("I94").Value2 = 1168.9333
$ws.Range("K94").Value2 = 1168.9333
$ws.Range("M94").Value2 = -717.9332999999999
$ws.Range("H98").Value2 = 1212.4348
$ws.Range("I98").Value2 = 1231.1818
$ws.Range("J98").Value2 = 800
$ws.Range("K98").Value2 = 1231.1818
$ws.Range("L98").Value2 = 800
$ws.Range("M98").Value2 = 266.8181999999999
$ws.Range("N98").Value2 = -3796
$ws.Range("H100").Value2 = 8507.143
$ws.Range("J100").Value2 = 8577.777
$ws.Range("L100").Value2 = 8577.777
$ws.Range("N100").Value2 = -9659.777
$ws.Range("H101").Value2 = 1124.625
$ws.Range("I101").Value2 = 704.7778
$ws.Range("K101").Value2 = 2114.3334
$ws.Range("M101").Value2 = -492.3334
$ws.Range("H106").Value2 = 2761.8
$ws.Range("I106").Value2 = 2827.25
$ws.Range("K106").Value2 = 2827.25
$ws.Range("M106").Value2 = -2196.25
$ws.Range("H122").Value2 = 1212.4348
$ws.Range("I122").Value2 = 1231.1818
$ws.Range("J122").Value2 = 800
$ws.Range("K122").Value2 = 3693.5454
$ws.Range("L122").Value2 = 2400
$ws.Range("M122").Value2 = -1243.5454
$ws.Range("N122").Value2 = -7300
$ws.Range("H132").Value2 = 15274.032
$ws.Range("I132").Value2 = 7692
$ws.Range("K132").Value2 = 23076
$ws.Range("M132").Value2 = -20546
$ws.Range("H135").Value2 = 11368136
$ws.Range("I135").Value2 = 14286782
$ws.Range("K135").Value2 = 128581038
$ws.Range("M135").Value2 = -128578503
$ws.Range("H137").Value2 = 2981.725
$ws.Range("I137").Value2 = 3097.1384
$ws.Range("J137").Value2 = 2481.6
$ws.Range("K137").Value2 = 9291.415199999999
$ws.Range("L137").Value2 = 7444.799999999999
$ws.Range("M137").Value2 = -6741.415199999999
$ws.Range("N137").Value2 = -12544.8
$ws.Range("H138").Value2 = 3063.7637
$ws.Range("I138").Value2 = 995.6539
$ws.Range("K138").Value2 = 2986.9617
$ws.Range("M138").Value2 = 2153.0383
$ws.Range("H141").Value2 = 1952.8306
$ws.Range("I141").Value2 = 1306.3773
$ws.Range("J141").Value2 = 7663.1665
$ws.Range("K141").Value2 = 3919.1319
$ws.Range("L141").Value2 = 22989.4995
$ws.Range("M141").Value2 = 1260.8681
$ws.Range("N141").Value2 = -33349.49950000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 3955.6162
$ws.Range("I32").Value2 = 3969.2842
$ws.Range("K32").Value2 = 3969.2842
$ws.Range("M32").Value2 = -3682.2842
$ws.Range("H50").Value2 = 4775
$ws.Range("J50").Value2 = 6350.6665
$ws.Range("L50").Value2 = 6350.6665
$ws.Range("N50").Value2 = -7778.6665
$ws.Range("H61").Value2 = 2330.2942
$ws.Range("I61").Value2 = 1124.3077
$ws.Range("J61").Value2 = 6249.75
$ws.Range("K61").Value2 = 1124.3077
$ws.Range("L61").Value2 = 6249.75
$ws.Range("M61").Value2 = -912.3077000000001
$ws.Range("N61").Value2 = -6673.75
$ws.Range("H74").Value2 = 2113.4375
$ws.Range("I74").Value2 = 1346
$ws.Range("K74").Value2 = 1346
$ws.Range("M74").Value2 = -472
$ws.Range("H77").Value2 = 2113.4375
$ws.Range("I77").Value2 = 1346
$ws.Range("K77").Value2 = 6730
$ws.Range("M77").Value2 = -2362
$ws.Range("H102").Value2 = 22224438
$ws.Range("I102").Value2 = 2172.2307
$ws.Range("K102").Value2 = 2172.2307
$ws.Range("M102").Value2 = -550.2307000000001
$ws.Range("H122").Value2 = 2534.4211
$ws.Range("I122").Value2 = 1314
$ws.Range("K122").Value2 = 3942
$ws.Range("M122").Value2 = -1492
$ws.Range("H132").Value2 = 1318.9166
$ws.Range("I132").Value2 = 1181.3
$ws.Range("K132").Value2 = 3543.9
$ws.Range("M132").Value2 = -1013.9
$ws.Range("H136").Value2 = 2330.2942
$ws.Range("I136").Value2 = 1124.3077
$ws.Range("J136").Value2 = 6249.75
$ws.Range("K136").Value2 = 3372.9231
$ws.Range("L136").Value2 = 18749.25
$ws.Range("M136").Value2 = -822.9231
$ws.Range("N136").Value2 = -23849.25
$ws.Range("H138").Value2 = 59824.668
$ws.Range("J138").Value2 = 59824.668
$ws.Range("L138").Value2 = 59824.668
$ws.Range("N138").Value2 = -70104.66800000001
$ws.Range("H139").Value2 = 80079.44500000001
$ws.Range("J139").Value2 = 80079.44500000001
$ws.Range("L139").Value2 = 80079.44500000001
$ws.Range("N139").Value2 = -90359.44500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value2 = 29747.666
$ws.Range("J58").Value2 = 18750
$ws.Range("L58").Value2 = 18750
$ws.Range("N58").Value2 = -19338
$ws.Range("H60").Value2 = 0
$ws.Range("J60").Value2 = 0
$ws.Range("L60").Value2 = 0
$ws.Range("N60").ClearContents()
$ws.Range("H105").Value2 = 3475.6538
$ws.Range("I105").Value2 = 2211.9285
$ws.Range("K105").Value2 = 2211.9285
$ws.Range("M105").Value2 = -464.9285
$ws.Range("H134").Value2 = 3759.4783
$ws.Range("I134").Value2 = 3356.3125
$ws.Range("K134").Value2 = 10068.9375
$ws.Range("M134").Value2 = -7533.9375
$ws.Range("H138").Value2 = 75068.80499999999
$ws.Range("I138").Value2 = 59599
$ws.Range("J138").Value2 = 75687.60000000001
$ws.Range("K138").Value2 = 59599
$ws.Range("L138").Value2 = 75687.60000000001
$ws.Range("M138").Value2 = -54459
$ws.Range("N138").Value2 = -85967.60000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 1696.5714
$ws.Range("I31").Value2 = 1555.6046
$ws.Range("J31").Value2 = 2162.8462
$ws.Range("K31").Value2 = 1555.6046
$ws.Range("L31").Value2 = 2162.8462
$ws.Range("M31").Value2 = -1260.6046
$ws.Range("N31").Value2 = -2752.8462
$ws.Range("H34").Value2 = 1696.5714
$ws.Range("I34").Value2 = 1555.6046
$ws.Range("J34").Value2 = 2162.8462
$ws.Range("K34").Value2 = 1555.6046
$ws.Range("L34").Value2 = 2162.8462
$ws.Range("M34").Value2 = -1353.6046
$ws.Range("N34").Value2 = -2566.8462
$ws.Range("H47").Value2 = 0
$ws.Range("I47").Value2 = 0
$ws.Range("K47").Value2 = 0
$ws.Range("M47").ClearContents()
$ws.Range("H58").Value2 = 2461.037
$ws.Range("I58").Value2 = 1811.75
$ws.Range("J58").Value2 = 4316.143
$ws.Range("K58").Value2 = 1811.75
$ws.Range("L58").Value2 = 4316.143
$ws.Range("M58").Value2 = -1608.75
$ws.Range("N58").Value2 = -4722.143
$ws.Range("H62").Value2 = 13000
$ws.Range("I62").Value2 = 13000
$ws.Range("K62").Value2 = 13000
$ws.Range("M62").Value2 = -12376
$ws.Range("H65").Value2 = 13000
$ws.Range("I65").Value2 = 13000
$ws.Range("K65").Value2 = 65000
$ws.Range("M65").Value2 = -61880
$ws.Range("H105").Value2 = 3642.8572
$ws.Range("I105").Value2 = 997.5
$ws.Range("K105").Value2 = 997.5
$ws.Range("M105").Value2 = 749.5
$ws.Range("H132").Value2 = 1981.2727
$ws.Range("I132").Value2 = 2108.1943
$ws.Range("J132").Value2 = 1410.125
$ws.Range("K132").Value2 = 6324.5829
$ws.Range("L132").Value2 = 4230.375
$ws.Range("M132").Value2 = -3794.5829
$ws.Range("N132").Value2 = -9290.375
$ws.Range("H134").Value2 = 3024.4443
$ws.Range("I134").Value2 = 2541.8718
$ws.Range("J134").Value2 = 4279.1333
$ws.Range("K134").Value2 = 7625.6154
$ws.Range("L134").Value2 = 12837.3999
$ws.Range("M134").Value2 = -5090.6154
$ws.Range("N134").Value2 = -17907.3999
$ws.Range("H136").Value2 = 2461.037
$ws.Range("I136").Value2 = 1811.75
$ws.Range("J136").Value2 = 4316.143
$ws.Range("K136").Value2 = 5435.25
$ws.Range("L136").Value2 = 12948.429
$ws.Range("M136").Value2 = -2885.25
$ws.Range("N136").Value2 = -18048.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value2 = 17649546
$ws.Range("I4").Value2 = 25795260
$ws.Range("K4").Value2 = 77385780
$ws.Range("M4").Value2 = -77385668
$ws.Range("H61").Value2 = 200.5
$ws.Range("I61").Value2 = 163
$ws.Range("J61").Value2 = 413
$ws.Range("K61").Value2 = 489
$ws.Range("L61").Value2 = 1239
$ws.Range("M61").Value2 = -274
$ws.Range("N61").Value2 = -1669
$ws.Range("H108").Value2 = 2427
$ws.Range("I108").Value2 = 2427
$ws.Range("K108").Value2 = 7281
$ws.Range("M108").Value2 = -4401
$ws.Range("H114").Value2 = 924.6667
$ws.Range("I114").Value2 = 799
$ws.Range("J114").Value2 = 987.5
$ws.Range("K114").Value2 = 2397
$ws.Range("L114").Value2 = 2962.5
$ws.Range("M114").Value2 = 857
$ws.Range("N114").Value2 = -9470.5
$ws.Range("H115").Value2 = 5000
$ws.Range("J115").Value2 = 5000
$ws.Range("L115").Value2 = 15000
$ws.Range("N115").Value2 = -17350
$ws.Range("H117").Value2 = 1960.5714
$ws.Range("I117").Value2 = 612
$ws.Range("J117").Value2 = 2500
$ws.Range("K117").Value2 = 1836
$ws.Range("L117").Value2 = 7500
$ws.Range("M117").Value2 = 1606
$ws.Range("N117").Value2 = -14384
$ws.Range("H121").Value2 = 6945342.5
$ws.Range("I121").Value2 = 291
$ws.Range("K121").Value2 = 873
$ws.Range("M121").Value2 = 437
$ws.Range("H132").Value2 = 1145.3334
$ws.Range("J132").Value2 = 1300
$ws.Range("L132").Value2 = 11700
$ws.Range("N132").Value2 = -16760

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value2 = 143072.72
$ws.Range("I3").Value2 = 143072.72
$ws.Range("K3").Value2 = 143072.72
$ws.Range("M3").Value2 = -142956.72
$ws.Range("H40").Value2 = 14993
$ws.Range("J40").Value2 = 14993
$ws.Range("L40").Value2 = 14993
$ws.Range("N40").Value2 = -15295
$ws.Range("H70").Value2 = 121671.2
$ws.Range("J70").Value2 = 12249.167
$ws.Range("L70").Value2 = 12249.167
$ws.Range("N70").Value2 = -12789.167
$ws.Range("H73").Value2 = 121671.2
$ws.Range("J73").Value2 = 12249.167
$ws.Range("L73").Value2 = 12249.167
$ws.Range("N73").Value2 = -14121.167
$ws.Range("H80").Value2 = 18582952
$ws.Range("I80").Value2 = 141225.75
$ws.Range("K80").Value2 = 141225.75
$ws.Range("M80").Value2 = -140227.75
$ws.Range("H83").Value2 = 18582952
$ws.Range("I83").Value2 = 141225.75
$ws.Range("K83").Value2 = 706128.75
$ws.Range("M83").Value2 = -701136.75
$ws.Range("H102").Value2 = 2758.7058
$ws.Range("I102").Value2 = 1684.6154
$ws.Range("K102").Value2 = 1684.6154
$ws.Range("M102").Value2 = -62.61539999999991
$ws.Range("H113").Value2 = 8843.777
$ws.Range("I113").Value2 = 5148.5
$ws.Range("J113").Value2 = 11800
$ws.Range("K113").Value2 = 5148.5
$ws.Range("L113").Value2 = 11800
$ws.Range("M113").Value2 = -2978.5
$ws.Range("N113").Value2 = -16140
$ws.Range("H122").Value2 = 6168.8335
$ws.Range("I122").Value2 = 4320.6665
$ws.Range("K122").Value2 = 12961.9995
$ws.Range("M122").Value2 = -10511.9995
$ws.Range("H132").Value2 = 3896.6667
$ws.Range("I132").Value2 = 4120.8335
$ws.Range("J132").Value2 = 3000
$ws.Range("K132").Value2 = 12362.5005
$ws.Range("L132").Value2 = 9000
$ws.Range("M132").Value2 = -9832.500499999998
$ws.Range("N132").Value2 = -14060
$ws.Range("H139").Value2 = 40000
$ws.Range("J139").Value2 = 40000
$ws.Range("L139").Value2 = 40000
$ws.Range("N139").Value2 = -50280
$ws.Range("H141").Value2 = 46000
$ws.Range("I141").Value2 = 0
$ws.Range("J141").Value2 = 46000
$ws.Range("K141").Value2 = 0
$ws.Range("L141").Value2 = 46000
$ws.Range("N141").Value2 = -56360
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value2 = 10399.941
$ws.Range("I46").Value2 = 3799
$ws.Range("J46").Value2 = 10812.5
$ws.Range("K46").Value2 = 3799
$ws.Range("L46").Value2 = 10812.5
$ws.Range("M46").Value2 = -3611
$ws.Range("N46").Value2 = -11188.5
$ws.Range("H55").Value2 = 1689.5294
$ws.Range("I55").Value2 = 1995
$ws.Range("J55").Value2 = 696.75
$ws.Range("K55").Value2 = 1995
$ws.Range("L55").Value2 = 696.75
$ws.Range("M55").Value2 = -1822
$ws.Range("N55").Value2 = -1042.75
$ws.Range("H61").Value2 = 1012.3333
$ws.Range("I61").Value2 = 951.375
$ws.Range("K61").Value2 = 951.375
$ws.Range("M61").Value2 = -749.375
$ws.Range("H100").Value2 = 98924.586
$ws.Range("I100").Value2 = 226219
$ws.Range("K100").Value2 = 226219
$ws.Range("M100").Value2 = -225678
$ws.Range("H113").Value2 = 1012.3333
$ws.Range("I113").Value2 = 951.375
$ws.Range("K113").Value2 = 951.375
$ws.Range("M113").Value2 = 1218.625
$ws.Range("H122").Value2 = 6770.7646
$ws.Range("I122").Value2 = 2804
$ws.Range("K122").Value2 = 8412
$ws.Range("M122").Value2 = -5962
$ws.Range("H132").Value2 = 3718.554
$ws.Range("I132").Value2 = 2811.8833
$ws.Range("K132").Value2 = 8435.6499
$ws.Range("M132").Value2 = -5905.6499
$ws.Range("H134").Value2 = 64000
$ws.Range("J134").Value2 = 64000
$ws.Range("L134").Value2 = 64000
$ws.Range("N134").Value2 = -74140
$ws.Range("H136").Value2 = 4118
$ws.Range("I136").Value2 = 4013.7058
$ws.Range("J136").Value2 = 5300
$ws.Range("K136").Value2 = 12041.1174
$ws.Range("L136").Value2 = 15900
$ws.Range("M136").Value2 = -9491.117400000001
$ws.Range("N136").Value2 = -21000
$ws.Range("H138").Value2 = 71497.5
$ws.Range("J138").Value2 = 71497.5
$ws.Range("L138").Value2 = 71497.5
$ws.Range("N138").Value2 = -81777.5
$ws.Range("H139").Value2 = 95523.89
$ws.Range("J139").Value2 = 95523.89
$ws.Range("L139").Value2 = 95523.89
$ws.Range("N139").Value2 = -105803.89
$ws.Range("H140").Value2 = 171388
$ws.Range("J140").Value2 = 171388
$ws.Range("L140").Value2 = 171388
$ws.Range("N140").Value2 = -181748
$ws.Range("H141").Value2 = 68900
$ws.Range("J141").Value2 = 66555.55499999999
$ws.Range("L141").Value2 = 66555.55499999999
$ws.Range("N141").Value2 = -76915.55499999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value2 = 2290.55
$ws.Range("I14").Value2 = 1446
$ws.Range("K14").Value2 = 1446
$ws.Range("M14").Value2 = -1278
$ws.Range("H107").Value2 = 471.28125
$ws.Range("I107").Value2 = 555.36365
$ws.Range("J107").Value2 = 286.3
$ws.Range("K107").Value2 = 1666.09095
$ws.Range("L107").Value2 = 858.9000000000001
$ws.Range("M107").Value2 = 253.90905
$ws.Range("N107").Value2 = -4698.9
$ws.Range("H113").Value2 = 687.0769
$ws.Range("I113").Value2 = 446
$ws.Range("J113").Value2 = 1699.6
$ws.Range("K113").Value2 = 1338
$ws.Range("L113").Value2 = 5098.799999999999
$ws.Range("M113").Value2 = 832
$ws.Range("N113").Value2 = -9438.799999999999
$ws.Range("H122").Value2 = 283518.16
$ws.Range("I122").Value2 = 2545.5518
$ws.Range("J122").Value2 = 1024264.06
$ws.Range("K122").Value2 = 7636.655400000001
$ws.Range("L122").Value2 = 3072792.18
$ws.Range("M122").Value2 = -5186.655400000001
$ws.Range("N122").Value2 = -3077692.18
$ws.Range("H123").Value2 = 52500
$ws.Range("J123").Value2 = 52500
$ws.Range("L123").Value2 = 52500
$ws.Range("N123").Value2 = -62300
$ws.Range("H132").Value2 = 1455.3864
$ws.Range("I132").Value2 = 1140.2433
$ws.Range("K132").Value2 = 3420.7299
$ws.Range("M132").Value2 = -890.7299000000003
$ws.Range("H135").Value2 = 77333.22
$ws.Range("J135").Value2 = 77333.22
$ws.Range("L135").Value2 = 77333.22
$ws.Range("N135").Value2 = -87473.22
$ws.Range("H136").Value2 = 1467.2593
$ws.Range("I136").Value2 = 809.35
$ws.Range("K136").Value2 = 2428.05
$ws.Range("M136").Value2 = 121.9499999999998
